$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plano de Ação")

# --- Fix existing label text in row 20 (was "Projeto do site", now "Tela Sobre nós") ---
$ws.Range("A20").Value = "Tela Sobre nós"

# --- Populate the "Ação" column (F) for the Sprint 3 block (rows 16-21), which was empty ---
# Match formatting of column E in each row first, then set the text.
$ws.Range("E16:E21").Copy()
$ws.Range("F16:F21").PasteSpecial(-4122)

$ws.Range("F16").Value = "Criar o CSS e o HTML do Detalhes do projeto"
$ws.Range("F17").Value = "Criar o css e html do footer e implementar no site"
$ws.Range("F18").Value = "Fazer o CSS e o HTML da Dashboard"
$ws.Range("F19").Value = "Criar o Script de validação das inputs e e mensagem de erro "
$ws.Range("F20").Value = "Fazer o CSS e o HTML da sobre nós"
$ws.Range("F21").Value = "Modificar o codigo da API para que ele interaja com o banco de dados"

# Rows grew taller once the Ação text wraps onto multiple lines.
$ws.Rows("16").RowHeight = 38.25
$ws.Rows("17").RowHeight = 27
$ws.Rows("18").RowHeight = 24
$ws.Rows("19").RowHeight = 40.5
$ws.Rows("20").RowHeight = 26.25
$ws.Rows("21").RowHeight = 46.5

# --- Add a new "Sprint 4 / Scrum Master SP2D: Matteus" block below the existing ones ---
# Duplicate the prior block's section header + column header + 6 data rows (A14:G21),
# then edit the pasted copy's content in place.
$ws.Range("A14:G21").Copy()
$ws.Range("A22").PasteSpecial()
$excel.CutCopyMode = 0

$ws.Range("A22").Value = "Scrum Master SP2D: Matteus"

$ws.Range("A24").Value = "Definir os graficos "
$ws.Range("D24").Value = 45222
$ws.Range("A25").Value = "Diagrama de solução"
$ws.Range("D25").Value = 45222
$ws.Range("A26").Value = "Escrever legenda no Diagrama de solução"
$ws.Range("D26").Value = 45222
$ws.Range("A27").Value = "Arrumar logica da calculadora"
$ws.Range("D27").Value = 45222
$ws.Range("A28").Value = "Modulo analytics "
$ws.Range("D28").Value = 45222

$ws.Rows("24").RowHeight = 26.25
$ws.Rows("25").RowHeight = 34.5
$ws.Rows("26").RowHeight = 34.5
$ws.Rows("27").RowHeight = 33
$ws.Rows("28").RowHeight = 28.5

# Row 29 stays as an empty spacer row, as in the rest of the sheet's blocks - clear
# the data that got copied into it from row 21.
$ws.Range("A29:G29").ClearContents()

# --- View bookkeeping: the author had scrolled/zoomed to the newly added rows ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 81
$ws.Range("F27").Select()
